$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 'Software Engineer'
$ws.Cells.Item(2, 2).Value = 'Advita Ortho'
$ws.Cells.Item(2, 3).Value = 'Gainesville, FL, US USA'
$ws.Cells.Item(2, 4).Value = 17.8
$ws.Cells.Item(2, 5).Value = 'RAG, Gemini, Hugging Face, Pinecone, ChromaDB, Prompt Engineering, TensorFlow, PyTorch, AWS SageMaker, Azure ML'
$ws.Cells.Item(2, 6).NumberFormat = "@"
$ws.Cells.Item(2, 6).Value = '2026-02-26'
$ws.Cells.Item(2, 6).Style = "Normal"
$ws.Cells.Item(2, 7).Value = 'https://www.indeed.com/viewjob?jk=88f54ab0bf4b2002'

# Row 3
$ws.Cells.Item(3, 1).Value = 'Senior Software Engineer, LLMs'
$ws.Cells.Item(3, 2).Value = 'Upstart'
$ws.Cells.Item(3, 3).Value = 'Remote, US USA'
$ws.Cells.Item(3, 4).Value = 16.7
$ws.Cells.Item(3, 5).Value = 'Generative AI, LangChain, RAG, LLaMA, FAISS, Pinecone, Prompt Engineering, FastAPI, Docker, Kubernetes'
$ws.Cells.Item(3, 6).NumberFormat = "@"
$ws.Cells.Item(3, 6).Value = '2026-02-26'
$ws.Cells.Item(3, 6).Style = "Normal"
$ws.Cells.Item(3, 7).Value = 'https://www.indeed.com/viewjob?jk=948dde4732205d6f'

# Row 4
$ws.Cells.Item(4, 1).Value = 'AI DevOps Engineer'
$ws.Cells.Item(4, 2).Value = 'Reef Capital Partners'
$ws.Cells.Item(4, 3).Value = 'Lehi, UT, US USA'
$ws.Cells.Item(4, 4).Value = 15.6
$ws.Cells.Item(4, 5).Value = 'LangChain, RAG, S3, Data Lake, AKS, CI/CD, GitHub Actions, Terraform, Git, Snowflake'
$ws.Cells.Item(4, 6).NumberFormat = "@"
$ws.Cells.Item(4, 6).Value = '2026-02-26'
$ws.Cells.Item(4, 6).Style = "Normal"
$ws.Cells.Item(4, 7).Value = 'https://www.indeed.com/viewjob?jk=1e695ed23053b341'

# Row 5
$ws.Cells.Item(5, 1).Value = 'Senior Python Backend Engineer'
$ws.Cells.Item(5, 2).Value = 'Sustainability Engineering Group'
$ws.Cells.Item(5, 3).Value = 'Phoenix, AZ, US USA'
$ws.Cells.Item(5, 4).Value = 15.6
$ws.Cells.Item(5, 5).Value = 'LangChain, RAG, LLaMA, Pinecone, FastAPI, Docker, CI/CD, GitHub Actions, Git, PostgreSQL'
$ws.Cells.Item(5, 6).NumberFormat = "@"
$ws.Cells.Item(5, 6).Value = '2026-02-26'
$ws.Cells.Item(5, 6).Style = "Normal"
$ws.Cells.Item(5, 7).Value = 'https://www.indeed.com/viewjob?jk=44862ba6d8cffad9'

# Row 6
$ws.Cells.Item(6, 1).Value = 'Data Scientist II'
$ws.Cells.Item(6, 2).Value = 'Avathon'
$ws.Cells.Item(6, 3).Value = 'Pleasanton, CA, US USA'
$ws.Cells.Item(6, 4).Value = 14.4
$ws.Cells.Item(6, 5).Value = 'Data Scientist, Generative AI, RAG, Hugging Face, Prompt Engineering, TensorFlow, PyTorch, CI/CD, Git, Python'
$ws.Cells.Item(6, 6).NumberFormat = "@"
$ws.Cells.Item(6, 6).Value = '2026-02-26'
$ws.Cells.Item(6, 6).Style = "Normal"
$ws.Cells.Item(6, 7).Value = 'https://www.indeed.com/viewjob?jk=2b79f4ac6055475e'

# Row 7
$ws.Cells.Item(7, 1).Value = 'Data Scientist I'
$ws.Cells.Item(7, 2).Value = 'Avathon'
$ws.Cells.Item(7, 3).Value = 'Pleasanton, CA, US USA'
$ws.Cells.Item(7, 4).Value = 13.3
$ws.Cells.Item(7, 5).Value = 'Data Scientist, Generative AI, RAG, Hugging Face, Prompt Engineering, TensorFlow, PyTorch, Git, Python, R'
$ws.Cells.Item(7, 6).NumberFormat = "@"
$ws.Cells.Item(7, 6).Value = '2026-02-26'
$ws.Cells.Item(7, 6).Style = "Normal"
$ws.Cells.Item(7, 7).Value = 'https://www.indeed.com/viewjob?jk=81f56a55e57fefe0'

# Row 8
$ws.Cells.Item(8, 1).Value = 'Senior Software Engineer (Java Full stack)'
$ws.Cells.Item(8, 2).Value = 'Optum'
$ws.Cells.Item(8, 3).Value = 'Raleigh, NC, US USA'
$ws.Cells.Item(8, 4).Value = 11.1
$ws.Cells.Item(8, 5).Value = 'Docker, Kubernetes, CI/CD, Terraform, Git, Kafka, SQL, R, Java, Scala'
$ws.Cells.Item(8, 6).NumberFormat = "@"
$ws.Cells.Item(8, 6).Value = '2026-02-26'
$ws.Cells.Item(8, 6).Style = "Normal"
$ws.Cells.Item(8, 7).Value = 'https://www.indeed.com/viewjob?jk=3fb41469a6c6781d'

# Row 9
$ws.Cells.Item(9, 1).Value = 'Sr Systems Engineer HPC'
$ws.Cells.Item(9, 2).Value = 'Shell'
$ws.Cells.Item(9, 3).Value = 'Houston, TX, US USA'
$ws.Cells.Item(9, 4).Value = 11.1
$ws.Cells.Item(9, 5).Value = 'RAG, Kubernetes, Terraform, Git, MySQL, Python, SQL, R, Java, Scala'
$ws.Cells.Item(9, 6).NumberFormat = "@"
$ws.Cells.Item(9, 6).Value = '2026-02-26'
$ws.Cells.Item(9, 6).Style = "Normal"
$ws.Cells.Item(9, 7).Value = 'https://www.indeed.com/viewjob?jk=19f6be6779a46e20'

# Row 10
$ws.Cells.Item(10, 1).Value = 'Sr Systems Engineer HPC'
$ws.Cells.Item(10, 2).Value = 'Shell'
$ws.Cells.Item(10, 3).Value = 'Houston, TX, US USA'
$ws.Cells.Item(10, 4).Value = 11.1
$ws.Cells.Item(10, 5).Value = 'RAG, Kubernetes, Terraform, Git, MySQL, Python, SQL, R, Java, Scala'
$ws.Cells.Item(10, 6).NumberFormat = "@"
$ws.Cells.Item(10, 6).Value = '2026-02-26'
$ws.Cells.Item(10, 6).Style = "Normal"
$ws.Cells.Item(10, 7).Value = 'https://www.indeed.com/viewjob?jk=23f4372fe821834a'

# Row 11
$ws.Cells.Item(11, 1).Value = 'Machine Learning Engineer, GenAI Technology'
$ws.Cells.Item(11, 2).Value = 'Point72'
$ws.Cells.Item(11, 3).Value = 'New York, NY, US USA'
$ws.Cells.Item(11, 4).Value = 10
$ws.Cells.Item(11, 5).Value = 'Data Scientist, Machine Learning Engineer, RAG, TensorFlow, PyTorch, Python, R, Java, Scala'
$ws.Cells.Item(11, 6).NumberFormat = "@"
$ws.Cells.Item(11, 6).Value = '2026-02-26'
$ws.Cells.Item(11, 6).Style = "Normal"
$ws.Cells.Item(11, 7).Value = 'https://www.indeed.com/viewjob?jk=cdf2b548f7a10938'

# Row 12
$ws.Cells.Item(12, 1).Value = 'Data Engineer I (Databricks, MLOps)'
$ws.Cells.Item(12, 2).Value = 'Travelers'
$ws.Cells.Item(12, 3).Value = 'Hartford, CT, US USA'
$ws.Cells.Item(12, 4).Value = 10
$ws.Cells.Item(12, 5).Value = 'Generative AI, RAG, MLflow, CI/CD, GitHub Actions, Terraform, Git, Databricks, R'
$ws.Cells.Item(12, 6).NumberFormat = "@"
$ws.Cells.Item(12, 6).Value = '2026-02-26'
$ws.Cells.Item(12, 6).Style = "Normal"
$ws.Cells.Item(12, 7).Value = 'https://www.indeed.com/viewjob?jk=a78961dd1a5a17bd'

# Row 13
$ws.Cells.Item(13, 1).Value = 'AI Search Engineer'
$ws.Cells.Item(13, 2).Value = 'NTT DATA'
$ws.Cells.Item(13, 3).Value = 'Plano, TX, US USA'
$ws.Cells.Item(13, 4).Value = 10
$ws.Cells.Item(13, 5).Value = 'FastAPI, Docker, Kubernetes, AKS, CI/CD, Git, Python, R, Scala'
$ws.Cells.Item(13, 6).NumberFormat = "@"
$ws.Cells.Item(13, 6).Value = '2026-02-18'
$ws.Cells.Item(13, 6).Style = "Normal"
$ws.Cells.Item(13, 7).Value = 'https://www.indeed.com/viewjob?jk=c988bd7e23d852be'

# Row 14
$ws.Cells.Item(14, 1).Value = 'Machine Learning Engineer'
$ws.Cells.Item(14, 2).Value = 'DocuSign'
$ws.Cells.Item(14, 3).Value = 'Seattle, WA, US USA'
$ws.Cells.Item(14, 4).Value = 10
$ws.Cells.Item(14, 5).Value = 'Machine Learning Engineer, Generative AI, Docker, Kubernetes, CI/CD, Git, Python, SQL, R'
$ws.Cells.Item(14, 6).NumberFormat = "@"
$ws.Cells.Item(14, 6).Value = '2026-02-26'
$ws.Cells.Item(14, 6).Style = "Normal"
$ws.Cells.Item(14, 7).Value = 'https://www.indeed.com/viewjob?jk=ce91aa193829c9fd'

# Row 15 is removed entirely in the updated match list
$ws.Rows("15:15").Delete()
